# Update "想去人数" (Column F) values across the sheets of the
# 杭州-漫展信息 workbook to match the regenerated gh-pages data
# (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @(
    @{Row=2;  Val=269},
    @{Row=3;  Val=1061},
    @{Row=4;  Val=9436},
    @{Row=5;  Val=199},
    @{Row=7;  Val=6467},
    @{Row=9;  Val=77},
    @{Row=10; Val=9896},
    @{Row=11; Val=11337},
    @{Row=13; Val=1167},
    @{Row=14; Val=4964},
    @{Row=16; Val=466},
    @{Row=18; Val=334},
    @{Row=20; Val=1345},
    @{Row=21; Val=261},
    @{Row=22; Val=1869},
    @{Row=23; Val=895},
    @{Row=24; Val=1270},
    @{Row=25; Val=859},
    @{Row=27; Val=2054},
    @{Row=28; Val=436},
    @{Row=29; Val=634},
    @{Row=30; Val=2697},
    @{Row=31; Val=188},
    @{Row=32; Val=1780},
    @{Row=33; Val=94},
    @{Row=34; Val=804},
    @{Row=35; Val=67},
    @{Row=36; Val=923},
    @{Row=37; Val=586},
    @{Row=38; Val=37},
    @{Row=39; Val=3362},
    @{Row=40; Val=239},
    @{Row=41; Val=87},
    @{Row=42; Val=519},
    @{Row=43; Val=586},
    @{Row=45; Val=900},
    @{Row=46; Val=246},
    @{Row=48; Val=4217},
    @{Row=49; Val=58}
)
foreach ($u in $sheet1Updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Val
}

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @(
    @{Row=9;  Val=12},
    @{Row=19; Val=0}
)
foreach ($u in $sheet2Updates) {
    $ws2.Cells.Item($u.Row, 6).Value = $u.Val
}

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @(
    @{Row=2; Val=5973}
)
foreach ($u in $sheet3Updates) {
    $ws3.Cells.Item($u.Row, 6).Value = $u.Val
}

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @(
    @{Row=2;  Val=269},
    @{Row=3;  Val=1061},
    @{Row=4;  Val=9436},
    @{Row=5;  Val=199},
    @{Row=8;  Val=6467},
    @{Row=10; Val=9896},
    @{Row=11; Val=11337},
    @{Row=13; Val=1167},
    @{Row=14; Val=4964},
    @{Row=16; Val=466},
    @{Row=18; Val=334},
    @{Row=21; Val=1345},
    @{Row=22; Val=261},
    @{Row=23; Val=1869},
    @{Row=24; Val=859},
    @{Row=26; Val=2054},
    @{Row=27; Val=436},
    @{Row=28; Val=634},
    @{Row=29; Val=2697},
    @{Row=30; Val=188},
    @{Row=31; Val=1780},
    @{Row=32; Val=94},
    @{Row=34; Val=804},
    @{Row=39; Val=67},
    @{Row=40; Val=923},
    @{Row=41; Val=586},
    @{Row=42; Val=37},
    @{Row=44; Val=239},
    @{Row=45; Val=586},
    @{Row=46; Val=900},
    @{Row=47; Val=246},
    @{Row=49; Val=4217}
)
foreach ($u in $sheet4Updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Val
}

$wb.Save()
